$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.899.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "'1.905.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'324.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "'0.4587"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("D8").Value = "'0.3816"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("D9").Value = "'0.07708"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("D10").Value = "'0.9786"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D11").Value = "'22.18"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.20%  "
$ws.Range("D12").Value = "'1.882.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.676"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.86%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'6.950"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "'0.07073"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "'83.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.43%  "
$ws.Range("D18").Value = "'0.000009445"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.00%  "
$ws.Range("D19").Value = "'16.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.05%  "
$ws.Range("D20").Value = "'0.9998"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "'28.894.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.28%  "
$ws.Range("D22").Value = "'5.309"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.87%  "
$ws.Range("D23").Value = "'10.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.16%  "
$ws.Range("D24").Value = "'2.098"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "'158.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "'19.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("D27").Value = "'5.664"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.91%  "
$ws.Range("D28").Value = "'117.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("D29").Value = "'1.877"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.28%  "
$ws.Range("D30").Value = "'0.09307"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").Value = "'0.8609"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").Value = "'5.088"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("D33").Value = "'1.242"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.52%  "
$ws.Range("D34").Value = "'3.030"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'0.05708"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("D36").Value = "'1.159"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").Value = "'1.001"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Value = "'0.02041"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.83%  "
$ws.Range("D39").Value = "'7.460"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.33%  "
$ws.Range("D40").Value = "'0.5487"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.99%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.919"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.72%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1754"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("D43").Value = "'9.325"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "'0.000002744"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -11.84%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'2.161"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.59%  "
$ws.Range("D46").Value = "'0.5173"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").Value = "'11.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("D48").Value = "'0.06889"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("D49").Value = "'1.778"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.76%  "
$ws.Range("D50").Value = "'110.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.70%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value = "'1.001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.10%  "
